{"js": "// Update the cover letter paragraph text. The paragraph is a single run\n// whose visual line breaks are <w:br/> elements; Office.js represents each\n// <w:br/> as a vertical-tab (\"\\v\" / \"\\u000B\") character inside Range/Paragraph\n// .text. We rebuild the whole paragraph text (5 \"blocks\" separated by a\n// blank line, i.e. two consecutive breaks) and replace it in one shot so the\n// run/break structure in the OOXML stays intact.\n\nconst blocks = [\n  \"I am a senior AI/ML Engineer with a decade of experience in the software development industry, specializing in AI, MLOps, and system architecture. My journey in AI began in America, where I have consistently been at the forefront of technological advancements, architecting and implementing cutting-edge software solutions tailored to business needs. My expertise in AI-driven platforms and projects has enabled me to streamline business operations, enhance customer service, and drive innovation.\",\n  \"At InsoftAI, I led the development of AI-driven platforms capable of handling up to 90% of customer inquiries, significantly enhancing operational efficiency. My work on Support-nGen\\u2122 and LLM Twin showcases my ability to develop proprietary systems that automate processes and generate creative ideas, facilitating brand creation and streamlining content creation. My proficiency in deploying scalable, secure, and efficient real-time predictions on AWS SageMaker further highlights my technical capabilities.\",\n  \"During my tenure at Brainhub, I developed Sierra.ai, revolutionizing document management and information accessibility, resulting in a 30% increase in operational efficiency for clients. My leadership in formulating technical strategies and optimizing multi-AI agents improved response times by 40% and reduced manual intervention. I have a proven track record of designing and implementing robust ML serving architectures and deploying scalable, cost-effective solutions that align with business goals.\",\n  \"My experience at Kensho involved building TTS and STT solutions, enhancing user experience in voice synthesis applications, and developing ML systems for forecasting energy consumption. I have demonstrated strong leadership abilities by mentoring junior staff and fostering skill development, enhancing team performance.\",\n  \"I am eager to define and drive the long-term ML technical strategy in alignment with product and business goals. My deep expertise in computer vision, GenAI, and adjacent fields, combined with my ability to lead and grow high-performing teams, makes me an ideal candidate for this role. I am committed to creating a team culture where people feel empowered, supported, and technically challenged, ensuring strong cross-functional collaboration and delivering state-of-the-art models into production swiftly. I look forward to contributing to your organization's success by leveraging my skills and experience in AI/ML engineering.\",\n];\n\n// Join the blocks with a pair of \"\\v\" (each \"\\v\" == one <w:br/>) so the\n// resulting text exactly reproduces \"<w:br/><w:br/>\" between blocks, same as\n// the original document.\nconst newText = blocks.join(\"\\v\\v\");\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst paragraph = paragraphs.items[0];\nparagraph.insertText(newText, \"Replace\");\nawait context.sync();\n", "ps1": "# Update the cover letter paragraph text. The paragraph is a single run\n# whose visual line breaks are <w:br/> elements; the Word COM object model\n# represents each <w:br/> as a vertical-tab (Chr(11)) character inside\n# Range.Text. We rebuild the whole paragraph text (5 \"blocks\" separated by a\n# blank line, i.e. two consecutive breaks) and assign it to the paragraph's\n# Range in one shot so the run/break structure in the OOXML stays intact.\n\n$trademark = [char]0x2122\n\n$blocks = @(\n    \"I am a senior AI/ML Engineer with a decade of experience in the software development industry, specializing in AI, MLOps, and system architecture. My journey in AI began in America, where I have consistently been at the forefront of technological advancements, architecting and implementing cutting-edge software solutions tailored to business needs. My expertise in AI-driven platforms and projects has enabled me to streamline business operations, enhance customer service, and drive innovation.\",\n    (\"At InsoftAI, I led the development of AI-driven platforms capable of handling up to 90% of customer inquiries, significantly enhancing operational efficiency. My work on Support-nGen\" + $trademark + \" and LLM Twin showcases my ability to develop proprietary systems that automate processes and generate creative ideas, facilitating brand creation and streamlining content creation. My proficiency in deploying scalable, secure, and efficient real-time predictions on AWS SageMaker further highlights my technical capabilities.\"),\n    \"During my tenure at Brainhub, I developed Sierra.ai, revolutionizing document management and information accessibility, resulting in a 30% increase in operational efficiency for clients. My leadership in formulating technical strategies and optimizing multi-AI agents improved response times by 40% and reduced manual intervention. I have a proven track record of designing and implementing robust ML serving architectures and deploying scalable, cost-effective solutions that align with business goals.\",\n    \"My experience at Kensho involved building TTS and STT solutions, enhancing user experience in voice synthesis applications, and developing ML systems for forecasting energy consumption. I have demonstrated strong leadership abilities by mentoring junior staff and fostering skill development, enhancing team performance.\",\n    \"I am eager to define and drive the long-term ML technical strategy in alignment with product and business goals. My deep expertise in computer vision, GenAI, and adjacent fields, combined with my ability to lead and grow high-performing teams, makes me an ideal candidate for this role. I am committed to creating a team culture where people feel empowered, supported, and technically challenged, ensuring strong cross-functional collaboration and delivering state-of-the-art models into production swiftly. I look forward to contributing to your organization's success by leveraging my skills and experience in AI/ML engineering.\"\n)\n\n# Join the blocks with a pair of Chr(11) (each Chr(11) == one <w:br/>) so the\n# resulting text exactly reproduces \"<w:br/><w:br/>\" between blocks, same as\n# the original document. Built manually (rather than [string]::Join) so the\n# control characters survive intact.\n$br = [char]11\n$newText = \"\"\nfor ($i = 0; $i -lt $blocks.Length; $i++) {\n    if ($i -gt 0) {\n        $newText += $br\n        $newText += $br\n    }\n    $newText += $blocks[$i]\n}\n\n$d = $word.ActiveDocument\n$paragraph = $d.Paragraphs.Item(1)\n$range = $paragraph.Range\n\n# Exclude the trailing paragraph mark so we only overwrite the visible text,\n# leaving the paragraph (and its end-of-paragraph run) intact.\n$range.End = $range.End - 1\n$range.Text = $newText\n"}
